$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the rows for years 2000-2009 (old rows 2-11), shifting 2010-2020 up.
$ws.Range("A2:A11").EntireRow.Delete() | Out-Null

# Add the new 2021 row at row 13.
$ws.Range("A13").Value = "2021年"
$ws.Range("B13").Value = 15894888.52
$ws.Range("C13").Value = 8127066.05
$ws.Range("D13").Value = 2941

# Copy style from the row above (2020, now row 12) to the new 2021 row (A13) to match formatting.
$ws.Range("A12").Copy() | Out-Null
$ws.Range("A13").PasteSpecial(-4122) | Out-Null
